# "Fruta / hortaliza, semanal" - add a new weekly price record for
# Ciboulette (Vega Central Mapocho de Santiago) and push the existing
# records down by one row, exactly like Excel's native "insert row above".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 350; everything that used
# to live at rows 350-434 shifts down to 351-435 (dimension grows to R435).
$ws.Rows(350).Insert()

# Populate the newly inserted row 350 with the new record's data.
$ws.Range('A350').Value = 9
$ws.Range('B350').Value = 'Vega Central Mapocho de Santiago'
$ws.Range('C350').Value = 'Metropolitana'
$ws.Range('D350').Value = 44782
$ws.Range('E350').Value = 13
$ws.Range('F350').Value = 100112039
$ws.Range('G350').Value = 'Ciboulette'
$ws.Range('H350').Value = 'Sin especificar'
$ws.Range('I350').Value = 'Segunda'
$ws.Range('J350').Value = 510
$ws.Range('K350').Value = 1500
$ws.Range('L350').Value = 1800
$ws.Range('M350').Value = 1647
$ws.Range('N350').Value = '$/docena de atados'
$ws.Range('O350').Value = 'Provincia de Chacabuco'
$ws.Range('P350').Value = 549
$ws.Range('Q350').Value = 3
$ws.Range('R350').Value = 'Hortaliza'
